$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing user's password value
$ws.Range("B2").Value = 8520

# Fill in the new "manager" user row, mirroring the formatting of row 2
$ws.Range("A2:I2").Copy()
$ws.Range("A3:I3").PasteSpecial(-4122) # xlPasteFormats
$ws.Rows("3:3").RowHeight = $ws.Rows("2:2").RowHeight

$ws.Range("A3").Value = "manager"
$ws.Range("B3").Value = 123
$ws.Range("C3").Value = "Boshqaruvchi"
$ws.Range("D3").Value = "Manager"
$ws.Range("E3").Value = "+998 91 6554321"
$ws.Range("F3").Value = "manager@voltstream.uz"
$ws.Range("G3").Value = "Farg'ona, O'zbekiston"
$ws.Range("H3").Value = 35998

$ws.Hyperlinks.Add($ws.Range("F3"), "manager@voltstream.uz")
$ws.Range("F2").Copy()
$ws.Range("F3").PasteSpecial(-4122) # xlPasteFormats

# Excel auto-creates a built-in "Hyperlink" cell style when adding the
# hyperlink above; it is unused once we restore F3's style from F2, so
# drop it to keep the style table tidy.
$wb.Styles.Item("Hyperlink").Delete()
